$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FF")

# New quarterly column of data is being added to the front of the
# existing quarter columns (D:K -> E:L), so insert a fresh column
# before column D and let Excel shift everything else one column right.
$ws.Range("D1").EntireColumn.Insert()

# The newly inserted column D has no per-cell formatting yet; copy the
# formatting from column E (which now holds what used to be column D)
# so the new column matches the existing look (date format in the
# header rows, number format in the data rows, etc.)
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# A handful of rows are pure labels / spacer rows with no cell at all
# in columns D:K, so undo the formatting paste there to keep them empty.
$emptyRows = @(5, 6, 36, 37, 78, 79)
foreach ($r in $emptyRows) {
    $ws.Cells.Item($r, 4).Clear()
}

# Populate the new column D with the latest quarter's figures.
$dValues = @{
    7 = 43373
    8 = 81400
    9 = 68800
    10 = 12600
    12 = 900
    13 = 0
    14 = 0
    15 = 0
    17 = 71200
    18 = 10200
    20 = 3300
    21 = 16200
    22 = 0
    23 = 13500
    24 = 4000
    25 = 0
    26 = 9400
    27 = 9400
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = -3300
    33 = 9400
    34 = 0
    35 = 9400
    38 = 43373
    41 = 199300
    42 = 103100
    43 = 22000
    44 = 48000
    45 = 1300
    46 = 373700
    47 = 0
    48 = 104300
    49 = 1400
    50 = 0
    51 = 0
    52 = 3800
    53 = 0
    54 = 483200
    57 = 34500
    58 = 0
    59 = 17600
    60 = 52100
    61 = 0
    62 = 33200
    63 = 0
    64 = 0
    65 = 0
    66 = 85400
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 115600
    73 = 0
    74 = 0
    75 = 0
    76 = 397800
    77 = 0
    80 = 43373
    81 = 9400
    83 = 2700
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 13800
    91 = -2000
    92 = 0
    93 = 0
    94 = -11600
    96 = -2600
    97 = 0
    98 = 0
    99 = 0
    100 = -2800
    101 = 0
    102 = -600
}

foreach ($r in $dValues.Keys) {
    $ws.Cells.Item($r, 4).Value = $dValues[$r]
}

Write-Host "Applied quarterly column insert + new figures"
